# Write test cases for admin login functionality in testcases.xlsx
# (mirrors commit: "write test cases for admin login functionality in testcases.xlsx")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 10 new rows after the existing 8 data rows (row 9..18) ---
$ws.Range("A9:A18").EntireRow.Insert()

# --- 2. Carry the standard bordered/wrap-text formatting down into the new rows ---
$ws.Range("A2:G2").Copy()
$ws.Range("A9:G18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Full target content for data rows 2..18 (S/N, Functionality/Description,
#         Test Inputs, Test Procedure, Expected Results, Actual Results) ---
$rows = @(
  @{ Row=2; SN=1; C="Validate that student with correct credentials can login"; D="Username: amy.ng.2009`nPassword: qwerty128"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login success, student home page displayed with correct bidding summary"; G="Login successful but redirected to 'plan bid' instead of 'home' page"; GStyle=3 },
  @{ Row=3; SN=2; C="Validate that student who did not key in username would not be able to login successfully and show error message"; D="Username: `nPassword: qwerty128"; E="Select student as their role, do not key in anything in the username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=4; SN=3; C="Validate that student who did not key in password would not be able to login successfully and show error message"; D="Username: amy.ng.2009`nPassword: "; E="Select student as their role, put username into username input field, do not key in anything in the password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=5; SN=4; C="Validate that student with incorrect username would not be able to login and show error message"; D="Username: amy.ng`nPassword: qwerty128"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=6; SN=5; C="Validate that student with incorrect password would not be able to login and show error message"; D="Username: amy.ng.2009`nPassword: 123456"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=7; SN=6; C="Validate that student with correct credentials can login (double check)"; D="Username: ben.ng.2009`nPassword: qwerty129"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login success, student home page displayed with correct bidding summary"; G="Login successful but wrong name was displayed on home page"; GStyle=3 },
  @{ Row=8; SN=7; C="Validate that student with correct credentials can login (double double check)"; D="Username: calvin.ng.2009`nPassword: qwerty130"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login success, student home page displayed with correct bidding summary"; G="Login successful, bidding summary matches expected result"; GStyle=4 },
  @{ Row=9; SN=8; C="Validate that admin with correct credentials can login"; D="Username: admin`nPassword: adminpassword"; E="Select admin as their role, put username into username input field, put password into password input field. Submit"; F="Login success, bidding details for admin is displayed."; G="Login unsuccessful and error displayed (got directed to student home page and hence causing error as admin do not have correct bidding summary)"; GStyle=3 },
  @{ Row=10; SN=9; C="Validate that admin with correct credentials can login (double check)"; D="Username: admin`nPassword: adminpassword"; E="Select admin as their role, put username into username input field, put password into password input field. Submit"; F="Login success, bidding details for admin is displayed."; G="Login successful, bidding admin page matches expected result"; GStyle=4 },
  @{ Row=11; SN=10; C="Validate that if admin did not key in username, they would not be able to login successfully and show error message"; D="Username: `nPassword: adminpassword"; E="Select admin as their role, do not key in anything in the username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=12; SN=11; C="Validate that if admin did not key in password, they would not be able to login successfully and show error message"; D="Username: admin`nPassword: "; E="Select admin as their role, put username into username input field, do not key in anything in the password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=13; SN=12; C="Validate that admin with incorrect username would not be able to login and show error message"; D="Username: admin123`nPassword: password"; E="Select admin as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=14; SN=13; C="Validate that admin with incorrect password would not be able to login and show error message"; D="Username: admin`nPassword: woshiadmin"; E="Select admin as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=15; SN=14; C="Ensure that admin would not be able to login if they key in their password in CAPS "; D="Username: admin`nPassword: ADMINPASSWORD"; E="Select admin as their role, put username into username input field, key in password in CAPS into password input field. Submit"; F="Login failed, show error message"; G="Login successful and was brought to the bidding admin page"; GStyle=3 },
  @{ Row=16; SN=15; C="Ensure that admin would not be able to login if they key in their password in CAPS - double check"; D="Username: admin`nPassword: ADMINPASSWORD"; E="Select admin as their role, put username into username input field, key in password in CAPS into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 },
  @{ Row=17; SN=16; C="Validate if student is able to login into admin home page using student's credentials"; D="Username: amy.ng.2009`nPassword: qwerty128"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Login successful and was brought to the bidding admin page"; GStyle=3 },
  @{ Row=18; SN=17; C="Validate if student is able to login into admin home page using student's credentials"; D="Username: amy.ng.2009`nPassword: qwerty128"; E="Select student as their role, put username into username input field, put password into password input field. Submit"; F="Login failed, show error message"; G="Matched expected results"; GStyle=4 }
)

foreach ($item in $rows) {
  $r = $item.Row
  $ws.Cells.Item($r, 1).Value = $item.SN
  $ws.Cells.Item($r, 2).Value = 1.1
  $ws.Cells.Item($r, 3).Value = $item.C
  $ws.Cells.Item($r, 4).Value = $item.D
  $ws.Cells.Item($r, 5).Value = $item.E
  $ws.Cells.Item($r, 6).Value = $item.F
  $ws.Cells.Item($r, 7).Value = $item.G

  # Colour-code the "Actual Results" column: red = unexpected/bug, green = matched expectation
  if ($item.GStyle -eq 3) {
    $ws.Cells.Item($r, 7).Interior.Color = 8487423
  } else {
    $ws.Cells.Item($r, 7).Interior.Color = 6280092
  }
}

Write-Host "Inserted $($rows.Count - 0) admin-login rows; wrote $($rows.Count) rows total (2..18)."

# --- 4. Leave the selection where the author ended up while reviewing the new rows ---
$ws.Range("E14").Select() | Out-Null
